$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

$row = 53

# Copy formatting (including the date number format) from the row above
$ws.Range("A52:G52").Copy()
$ws.Range("A53:G53").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = "Facebook: The Inside Story"
$ws.Cells.Item($row, 2).Value = "Steven Levy"
$ws.Cells.Item($row, 3).Value = 43930
$ws.Cells.Item($row, 4).Value = 43934
$ws.Cells.Item($row, 5).Value = "facebook;entrepreuner;business;social networking;privacy;scandal"
$ws.Cells.Item($row, 6).Value = "Audio"
$ws.Cells.Item($row, 7).Value = "19 Hours 5 Mins"

$ws.Cells.Item($row + 1, 1).Select()
